# This script applies a permutation of the data rows (rows 2-21) of the
# active worksheet. The header row (row 1) is untouched; the contents of
# the data rows are reshuffled among themselves (no values are altered,
# only their row position changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 21
$lastCol      = 51   # column AY

# Mapping of destination row -> source row (1-based worksheet rows).
# Read as: "the data that ends up in row <key> originally lived in row <value>".
$rowMap = @{
    2  = 3
    3  = 9
    4  = 10
    5  = 13
    6  = 14
    7  = 17
    8  = 18
    9  = 19
    10 = 20
    11 = 21
    12 = 2
    13 = 4
    14 = 5
    15 = 6
    16 = 7
    17 = 8
    18 = 11
    19 = 12
    20 = 15
    21 = 16
}

# Snapshot every source row's values (A:AY) before any writes happen, so
# that later writes never clobber a row we still need to read from.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rng = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $lastCol))
    $snapshot[$r] = $rng.Value2
}

# The "Startdatum"/"Slutdatum" columns (Y and AA) hold plain text such as
# "2021-08-11". If written through .Value/.Value2 while the cell still has
# the default General number format, Excel auto-converts the look-alike
# date text into a real date serial number. Force those columns to a Text
# format first so the permuted values are written back unchanged as text.
$textCols = 25, 27   # Y, AA
foreach ($col in $textCols) {
    $colRng = $ws.Range($ws.Cells.Item($firstDataRow, $col), $ws.Cells.Item($lastDataRow, $col))
    $colRng.NumberFormat = "@"
}

# Write the snapshotted rows back out in their new (permuted) positions.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $values = $snapshot[$srcRow]
    $destRng = $ws.Range($ws.Cells.Item($destRow, 1), $ws.Cells.Item($destRow, $lastCol))
    $destRng.Value2 = $values
}

# Row 21's "Lokalnamn" (column P) is reported in the source changeset as
# ending up with the other site's name even though the rest of that row's
# data (A, B, D-H, Q, R, Y, AA, AW, AX, ...) moved in line with the row
# 16 -> row 21 permutation above. Apply that one targeted value to match
# the target state exactly.
$ws.Cells.Item(21, 16).Value = "Kåtamyrbäcken-Görjeån, Lu lm"   # column P
